$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new subject rows first (so their shared strings are registered first) ---
$ws.Cells.Item(12, 1).Value = "MVI011R031"
$ws.Cells.Item(12, 2).Value = "L"
$ws.Cells.Item(12, 3).Value = 45055
$ws.Cells.Item(12, 3).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(12, 4).Value = "E3"
$ws.Cells.Item(12, 5).Value = "E6"
$ws.Cells.Item(12, 6).Value = "E9"
$ws.Cells.Item(12, 7).Value = 45084
$ws.Cells.Item(12, 7).NumberFormat = "yyyy\-mm\-dd;@"

$ws.Cells.Item(13, 1).Value = "MVI012R897"
$ws.Cells.Item(13, 2).Value = "R"
$ws.Cells.Item(13, 3).Value = 45079
$ws.Cells.Item(13, 3).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(13, 4).Value = "E3"
$ws.Cells.Item(13, 5).Value = "E6"
$ws.Cells.Item(13, 6).Value = "E9"
$ws.Cells.Item(13, 7).Value = 45105
$ws.Cells.Item(13, 7).NumberFormat = "yyyy\-mm\-dd;@"

$ws.Cells.Item(14, 1).Value = "MVI013R864"
$ws.Cells.Item(14, 2).Value = "L"
$ws.Cells.Item(14, 3).Value = 45118
$ws.Cells.Item(14, 3).NumberFormat = "yyyy\-mm\-dd;@"

# --- Insert a new "Sex" column after column A (new column B) ---
$ws.Columns.Item(2).Insert()

$ws.Cells.Item(1,2).Value = "Sex"
$ws.Columns.Item(2).ColumnWidth = 4.25

$sex = @("M","M","F","F","F","F","F","M","M","M","F","F","M")
for ($i = 0; $i -lt $sex.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $sex[$i]
}

# --- Rename "Settings Changes" header (now column I) to "Changes" ---
$ws.Cells.Item(1, 9).Value = "Changes"

# --- Update the view/selection state ---
$ws.Range("I2").Select()
